$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Sprint 4 data: bump the sprint count and add the new sprint's story points
$ws.Range("B5").Value = 4
$ws.Range("D8").Value = 13

# Update the average formula to include the new sprint value
$ws.Range("C16").Formula = "=(D5+D6+D7+D8)/B5"

# Move the active selection to reflect where the new entry was made
$ws.Range("D9").Select()

# Recalculate and extend the chart's source range to include the new sprint
$wb.Application.Calculate()

$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Values = "=Tabelle1!`$D`$5:`$D`$8"
